# Actualización automática 2025-05-30 16:20:08
# Applies the data refresh to the "VENTAS POR GRUPO" sheet:
#  - Narrows column I (9) to width 9
#  - Updates a batch of numeric sales figures
#  - Refreshes the "X de 54" summary counters in row 56

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Column width change (stored width 12 -> 9) ---
# Excel's ColumnWidth property is expressed in character units and gets
# padded/rounded when persisted as the sheet's stored "width" attribute,
# so 8.17 is the character-width value that round-trips to a stored
# width of exactly 9 for this sheet's default font.
$ws.Columns.Item(9).ColumnWidth = 8.17

# --- Numeric cell updates ---
$ws.Range("L2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("L3").Value = -10.76
$ws.Range("L4").Value = 0
$ws.Range("L6").Value = 0

$ws.Range("C10").Value = 1002.24
$ws.Range("D10").Value = 0
$ws.Range("L10").Value = 17247.09
$ws.Range("M10").Value = 0

$ws.Range("D12").Value = -407.04
$ws.Range("L12").Value = 2295.56

$ws.Range("D16").Value = 0
$ws.Range("L16").Value = 6711.78

$ws.Range("C26").Value = 4624.12
$ws.Range("K26").Value = 2509.05
$ws.Range("L26").Value = 1605.8

$ws.Range("L27").Value = 768.95
$ws.Range("N27").Value = 3096.23

$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 44.89
$ws.Range("K29").Value = 743.08

$ws.Range("L42").Value = 1167.85

$ws.Range("L43").Value = 737.72

$ws.Range("D46").Value = 518.98
$ws.Range("E46").Value = 79.59999999999999

$ws.Range("L47").Value = 2436.41
$ws.Range("N47").Value = 1034.04

$ws.Range("L48").Value = 0

$ws.Range("L49").Value = 0

$ws.Range("L55").Value = 808.39

# --- Row 56 summary counters ("N de 54") ---
$ws.Range("C56").Value = "2 de 54"
$ws.Range("D56").Value = "1 de 54"
$ws.Range("E56").Value = "1 de 54"
$ws.Range("F56").Value = "0 de 54"
$ws.Range("G56").Value = "0 de 54"
$ws.Range("H56").Value = "0 de 54"
$ws.Range("I56").Value = "0 de 54"
$ws.Range("J56").Value = "1 de 54"
$ws.Range("K56").Value = "2 de 54"
$ws.Range("L56").Value = "9 de 54"
$ws.Range("M56").Value = "0 de 54"
$ws.Range("N56").Value = "2 de 54"
